$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '27.154.41'
Set-TextValue $ws.Range("E2") '  -0.24%  '
Set-TextValue $ws.Range("D3") '1.892.67'
Set-TextValue $ws.Range("E3") '  -0.75%  '
Set-TextValue $ws.Range("E4") '  -0.06%  '
Set-TextValue $ws.Range("D5") '306.94'
Set-TextValue $ws.Range("E5") '  -0.43%  '
Set-TextValue $ws.Range("E6") '  +0.03%  '
Set-TextValue $ws.Range("D7") '0.5222'
Set-TextValue $ws.Range("E7") '  -0.41%  '
Set-TextValue $ws.Range("D8") '0.3761'
Set-TextValue $ws.Range("E8") '  -0.54%  '
Set-TextValue $ws.Range("D9") '0.07274'
Set-TextValue $ws.Range("E9") '  +0.02%  '
Set-TextValue $ws.Range("D10") '21.11'
Set-TextValue $ws.Range("E10") '  -0.72%  '
Set-TextValue $ws.Range("D11") '0.9001'
Set-TextValue $ws.Range("E11") '  +0.09%  '
Set-TextValue $ws.Range("D12") '0.08163'
Set-TextValue $ws.Range("E12") '  +6.26%  '
Set-TextValue $ws.Range("D13") '1.921.64'
Set-TextValue $ws.Range("E13") '  +0.66%  '
Set-TextValue $ws.Range("D14") '96.40'
Set-TextValue $ws.Range("E14") '  +1.21%  '
Set-TextValue $ws.Range("D15") '5.287'
Set-TextValue $ws.Range("E15") '  +0.25%  '
Set-TextValue $ws.Range("E16") '  +0.08%  '
Set-TextValue $ws.Range("D17") '0.000008578'
Set-TextValue $ws.Range("E17") '  -1.00%  '
Set-TextValue $ws.Range("D18") '14.57'
Set-TextValue $ws.Range("E18") '  +0.21%  '
Set-TextValue $ws.Range("D20") '27.161.98'
Set-TextValue $ws.Range("E20") '  -0.50%  '
Set-TextValue $ws.Range("D21") '5.083'
Set-TextValue $ws.Range("E21") '  -0.11%  '
Set-TextValue $ws.Range("D22") '10.71'
Set-TextValue $ws.Range("E22") '  +0.67%  '
Set-TextValue $ws.Range("D23") '6.403'
Set-TextValue $ws.Range("E23") '  -0.71%  '
Set-TextValue $ws.Range("D24") '147.67'
Set-TextValue $ws.Range("E24") '  +1.44%  '
Set-TextValue $ws.Range("E25") '  -1.55%  '
Set-TextValue $ws.Range("D26") '18.18'
Set-TextValue $ws.Range("E26") '  +0.17%  '
Set-TextValue $ws.Range("D27") '1.743'
Set-TextValue $ws.Range("E27") '  +0.34%  '
Set-TextValue $ws.Range("D28") '114.85'
Set-TextValue $ws.Range("E28") '  -0.02%  '
Set-TextValue $ws.Range("D29") '4.895'
Set-TextValue $ws.Range("E29") '  -1.50%  '
Set-TextValue $ws.Range("D30") '4.789'
Set-TextValue $ws.Range("E30") '  -0.55%  '
Set-TextValue $ws.Range("D31") '0.09221'
Set-TextValue $ws.Range("E31") '  -0.19%  '
Set-TextValue $ws.Range("D32") '0.05051'
Set-TextValue $ws.Range("E32") '  -0.49%  '
Set-TextValue $ws.Range("D33") '0.7885'
Set-TextValue $ws.Range("E33") '  -2.47%  '
Set-TextValue $ws.Range("D34") '1.216'
Set-TextValue $ws.Range("E34") '  -2.36%  '
Set-TextValue $ws.Range("E35") '  +3.24%  '
Set-TextValue $ws.Range("D36") '2.958'
Set-TextValue $ws.Range("E36") '  -1.27%  '
Set-TextValue $ws.Range("D37") '2.582'
Set-TextValue $ws.Range("E37") '  -0.83%  '
Set-TextValue $ws.Range("D38") '0.5683'
Set-TextValue $ws.Range("E38") '  +0.15%  '
Set-TextValue $ws.Range("D39") '0.01987'
Set-TextValue $ws.Range("E39") '  -0.41%  '
Set-TextValue $ws.Range("E40") '  -0.14%  '
Set-TextValue $ws.Range("D41") '9.018'
Set-TextValue $ws.Range("E41") '  +0.15%  '
Set-TextValue $ws.Range("D42") '6.561'
Set-TextValue $ws.Range("E42") '  -1.25%  '
Set-TextValue $ws.Range("D43") '116.17'
Set-TextValue $ws.Range("E43") '  -2.64%  '
Set-TextValue $ws.Range("D44") '0.1519'
Set-TextValue $ws.Range("E44") '  -0.06%  '
Set-TextValue $ws.Range("D45") '0.4866'
Set-TextValue $ws.Range("E45") '  +0.42%  '
Set-TextValue $ws.Range("E46") '  +0.11%  '
Set-TextValue $ws.Range("D47") '10.06'
Set-TextValue $ws.Range("E47") '  -2.05%  '
Set-TextValue $ws.Range("D48") '1.623'
Set-TextValue $ws.Range("E48") '  +0.32%  '
Set-TextValue $ws.Range("D49") '38.09'
Set-TextValue $ws.Range("E49") '  +1.44%  '
Set-TextValue $ws.Range("D50") '63.40'
Set-TextValue $ws.Range("E50") '  -0.84%  '
Set-TextValue $ws.Range("E51") '  -0.03%  '
